# updated the data for 15 and 16th may
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Marksheet")

# Row 25 - Biology, 15 May 2025: no_of_questions=50, correct=46, incorrect=0, unattempted=4
$ws.Range("D25").Value = 50
$ws.Range("E25").Value = 46
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 4

# Row 27 - Chemistry, 16 May 2025: no_of_questions=25, correct=23, incorrect=0, unattempted=2
$ws.Range("D27").Value = 25
$ws.Range("E27").Value = 23
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 2

# Row 28 - Biology, 16 May 2025: no_of_questions=50, correct=49, incorrect=0, unattempted=1
$ws.Range("D28").Value = 50
$ws.Range("E28").Value = 49
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 1

# Match the author's final selection/cursor position in the sheet.
$ws.Range("E24").Select()
